$wb = $excel.ActiveWorkbook

$data = @{
  "2025" = @{ A2 = 0;                  B2 = 290.0628494009815;  E2 = 29049.07128553875;  G2 = 8095.925712662028;  I2 = 14940.02916486277;  L2 = 50998.86069102;     M2 = 11228.70813999;   N2 = 7234.065805482215; O2 = 6709.289187938176 }
  "2030" = @{ A2 = 219.6192975232897;  B2 = 3803.736742006061;  E2 = 45497.54827019678;  G2 = 8095.925712662028;  I2 = 31289.53350734562;  L2 = 60054.94214326091;  M2 = 17372.009741075;  N2 = 9195.86092813068;  O2 = 7881.440454293026 }
  "2035" = @{ A2 = 2152.642197777731;  B2 = 5760.108587711311;  E2 = 57498.34500251195;  G2 = 8095.925712662028;  I2 = 48287.37981995431;  L2 = 60054.94214326091;  M2 = 23337.84238116651; N2 = 13671.28110402295; O2 = 13096.12162563002 }
  "2040" = @{ A2 = 2152.642197777731;  B2 = 5760.108587711311;  E2 = 57498.34500251195;  G2 = 8095.925712662028;  I2 = 48287.37981995431;  L2 = 60054.94214326091;  M2 = 23337.84238116651; N2 = 13671.28110402295; O2 = 13096.12162563002 }
  "2045" = @{ A2 = 2152.642197777731;  B2 = 5760.108587711311;  E2 = 57498.34500251195;  G2 = 8095.925712662028;  I2 = 48287.37981995431;  L2 = 60054.94214326091;  M2 = 23337.84238116651; N2 = 13671.28110402295; O2 = 13096.12162563002 }
  "2050" = @{ A2 = 2152.642197777731;  B2 = 5760.108587711311;  E2 = 57498.34500251195;  G2 = 8095.925712662028;  I2 = 48287.37981995431;  L2 = 60054.94214326091;  M2 = 23337.84238116651; N2 = 13671.28110402295; O2 = 13096.12162563002 }
}

foreach ($sheetName in $data.Keys) {
  $ws = $wb.Worksheets.Item([string]$sheetName)
  $vals = $data[$sheetName]
  foreach ($cell in $vals.Keys) {
    $ws.Range([string]$cell).Value = $vals[$cell]
  }
}
